$wb = $excel.ActiveWorkbook

# --- Select I2 on "ResultCheckSheet2" (it currently holds tabSelected) ---
# Doing this now, while it is still the active sheet, records its last
# selection before we switch away to the newly added sheet.
$resultCheckSheet2 = $wb.Worksheets.Item("ResultCheckSheet2")
$resultCheckSheet2.Activate() | Out-Null
$resultCheckSheet2.Range("I2").Select() | Out-Null

# --- Add the new "EvaluateFormula" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "EvaluateFormula"

# --- Headers ---
$newSheet.Range("A1").Value = "Today"
$newSheet.Range("B1").Value = "Year"
$newSheet.Range("D1").Value = "Day"
$newSheet.Range("C1").Value = "Month"

# --- Re-use the existing date number-format style (m/d/yyyy) from the
#     "DataTypeAndFormatPattern" sheet so we don't create a duplicate style ---
$dataTypeSheet = $wb.Worksheets.Item("DataTypeAndFormatPattern")
$dataTypeSheet.Range("B2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- Formulas ---
$newSheet.Range("A2").Formula = "=TODAY()"
$newSheet.Range("B2").Formula = "=YEAR(A2)"
$newSheet.Range("C2").Formula = "=MONTH(A2)"
$newSheet.Range("D2").Formula = "=DAY(A2)"

# --- Page setup to match the rest of the workbook ---
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- Column A width (matches width="9" in the committed workbook) ---
$newSheet.Columns("A:A").ColumnWidth = 8.285714285714286

$newSheet.Range("A1").Select()
